$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.936.70"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "1.910.87"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.76%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9989"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5015"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3821"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07317"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9123"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07673"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("D13").Value = "1.920.28"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.491"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9988"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008747"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9986"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "27.976.84"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.186"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "2.169.44"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.612"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.841"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.80%  "
$ws.Range("E27").Value = "  +3.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.937"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09032"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.206"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.864"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.22%  "
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7798"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02087"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.606"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.068"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5554"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.05299"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.885"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "113.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.533"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4848"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9987"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.644"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06051"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.46%  "
